$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: swap Usinagem out of C3, put "-" ; put EAP into E3
$ws.Range("C3").Value = "-"
$ws.Range("E3").Value = "EAP"

# Row 6: swap EAP out of E6, put "-" ; put Usinagem into F6
$ws.Range("E6").Value = "-"
$ws.Range("F6").Value = "Usinagem"
